$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.471.43"
$ws.Range("D3").Value = "3.044.00"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.69"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.26"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.041.74"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  -6.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.50"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "3.554.91"
$ws.Range("D16").Value = "63.544.89"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "3.051.01"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.64"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.63"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.88"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.45"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.91"
$ws.Range("E24").Value = "  +7.93%  "
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.87"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.96"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.43"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.54"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.10"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.83"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0402"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.88"
$ws.Range("E38").Value = "  +12.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "434.09"
$ws.Range("E39").Value = "  -7.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0799"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").Value = "2.939.71"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.11"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.112"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.82"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.256"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.66"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "0.0₃0507"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("E51").Value = "  -1.89%  "
